$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$nl = [char]11
$table.Cell(1, 1).Range.Text = "80 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"
$table.Cell(1, 2).Range.Text = "23 x 62" + $nl + "  6    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "3|    |"
$table.Cell(1, 3).Range.Text = "79 x 91" + $nl + "  9    1" + $nl + "  ----" + $nl + "7|    |" + $nl + "9|    |"
$table.Cell(2, 1).Range.Text = "44 x 57" + $nl + "  5    7" + $nl + "  ----" + $nl + "4|    |" + $nl + "4|    |"
$table.Cell(2, 2).Range.Text = "50 x 99" + $nl + "  9    9" + $nl + "  ----" + $nl + "5|    |" + $nl + "0|    |"
$table.Cell(2, 3).Range.Text = "71 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$table.Cell(3, 1).Range.Text = "81 x 66" + $nl + "  6    6" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$table.Cell(3, 2).Range.Text = "30 x 64" + $nl + "  6    4" + $nl + "  ----" + $nl + "3|    |" + $nl + "0|    |"
$table.Cell(3, 3).Range.Text = "99 x 91" + $nl + "  9    1" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"
$table.Cell(4, 1).Range.Text = "89 x 79" + $nl + "  7    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "9|    |"
$table.Cell(4, 2).Range.Text = "18 x 32" + $nl + "  3    2" + $nl + "  ----" + $nl + "1|    |" + $nl + "8|    |"
$table.Cell(4, 3).Range.Text = "38 x 25" + $nl + "  2    5" + $nl + "  ----" + $nl + "3|    |" + $nl + "8|    |"
$table.Cell(5, 1).Range.Text = "96 x 73" + $nl + "  7    3" + $nl + "  ----" + $nl + "9|    |" + $nl + "6|    |"
$table.Cell(5, 2).Range.Text = "68 x 12" + $nl + "  1    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "8|    |"
$table.Cell(5, 3).Range.Text = "63 x 49" + $nl + "  4    9" + $nl + "  ----" + $nl + "6|    |" + $nl + "3|    |"
